# Remove the forecast blocks for Bắc Kinh, Thượng Hải and Tokyo from the
# "foreign_regression" sheet. Each city occupies 5 consecutive columns
# (Linear Regression, Prophet, 10% Growth, GDP Growth, Exponential Smoothing):
#
#   AP:AT = Bắc Kinh (Trung Quốc)
#   AU:AY = Thượng Hải (Trung Quốc)
#   AZ:BD = Thâm Quyến (Trung Quốc)   <- keep
#   BE:BI = Tokyo (Nhật Bản)
#   BJ:BN = Osaka (Nhật Bản)          <- keep
#   BO:BS = TP. Hồ Chí Minh (Việt Nam) <- keep
#
# Deleting these column ranges with a left-shift moves the remaining city
# blocks (Thâm Quyến, Osaka, TP. Hồ Chí Minh) into AP:BD, matching the
# target layout where the sheet's dimension shrinks from A1:BS13 to A1:BD13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from right to left so earlier column letters stay valid.
$ws.Range("BE1:BI13").Delete("ShiftLeft")
$ws.Range("AU1:AY13").Delete("ShiftLeft")
$ws.Range("AP1:AT13").Delete("ShiftLeft")
